# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.794.42"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.095.43"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'228.90"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "'61.58"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").Value = "'0.0846"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  +4.81%  "
$ws.Range("D13").Value = "2.405.88"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "'22.13"
$ws.Range("E15").Value = "  +4.61%  "
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "2.090.11"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "38.758.53"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").Value = "'71.93"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").Value = "'227.91"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").Value = "'171.39"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "'9.55"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("E29").Value = "  +3.59%  "
$ws.Range("D30").Value = "'19.35"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").Value = "'0.0619"
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("D36").Value = "'6.53"
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'3.60"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'18.15"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Value = "'0.0228"
$ws.Range("E41").Value = "  +4.22%  "
$ws.Range("D42").Value = "'101.33"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").Value = "1.535.81"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").Value = "'0.0910"
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'7.74"
$ws.Range("E46").Value = "  +6.57%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "'1.14"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").Value = "'4.12"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D51").Value = "2.289.53"
$ws.Range("E51").Value = "  -0.09%  "
